$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You have an application running on an EC2 instance that requires access to read/write to a DynamoDB table. Which of the following is a secure way for the EC2 instance to access the DynamoDB table?",
        "ques_type": 2,
        "options": [
            "Since the EC2 instance and the DynamoDB table are owned by the same AWS account, they can communicate with each other by default.",
            "Generate API keys with access to the DynamoDB table, and use the API keys in your application.",
            "Create an IAM user with the proper permissions, and associate the user with the EC2 instance.",
            "Create an IAM role with the proper permissions, and attach the role to the EC2 instance."
        ],
        "score": "Create an IAM role with the proper permissions, and attach the role to the EC2 instance."
    },
    {
        "title": "You currently have a mobile application that uploads files to an S3 bucket. When a new object is written to the S3 bucket, you need to store the metadata regarding the object into a DynamoDB table.Which of the following is the best solution to achieve this while minimizing your operational overhead?",
        "ques_type": 2,
        "options": [
            "Enable S3 event notifications to trigger a Lambda that will handle the logic.",
            "Enable S3 event notifications to directly trigger a write on your target DynamoDB table.",
            "Launch an EC2 instance that will host a script to periodically write metadata to the DynamoDB table.",
            "Configure a CloudWatch Events rule that will trigger a Lambda based on a cron expression."
        ],
        "score": "Enable S3 event notifications to trigger a Lambda that will handle the logic."
    },
    {
        "title": "Your company is currently hosting a MySQL database on-premise. Due to a lack of IT staff to handle the maintenance of the database, you\u2019ve been tasked to migrate it to AWS. You must ensure the migration will not require a major refactoring of the application code or a lot of time spent on maintenance in the future. Which of the following satisfies all the requirements?",
        "ques_type": 2,
        "options": [
            "Migrate the database to AWS DynamoDB.",
            "Migrate the database to AWS Redshift.",
            "Migrate the database to AWS RDS.",
            "Migrate the database to an EC2 instance built off of an Amazon-provided MySQL AMI."
        ],
        "score": "Migrate the database to AWS RDS."
    },
    {
        "title": "You have to store archival data that must be retrievable within a maximum of eight hours. In addition to the retrieval-time requirement, you must also minimize cost. What is the best solution to achieve this?",
        "ques_type": 2,
        "options": [
            "AWS S3 Standard storage",
            "AWS S3 Infrequent Access",
            "AWS Glacier",
            "AWS RDS"
        ],
        "score": "AWS Glacier"
    }
]
'@

# Update the questions text (currently in A2) to the reformatted JSON.
$ws.Range("A2").Value = $newText
$ws.Rows.Item(2).EntireRow.AutoFit()

# Delete row 1 (the old numeric placeholder cell, with its bold/bordered style).
# This shifts row 2 -- now holding the updated text -- up into row 1, and it
# keeps row 2's own (unstyled) formatting rather than inheriting row 1's style.
$ws.Rows.Item(1).Delete()
